# Insert a new "Quantite" column (numeric quantities) before the existing
# "Quantité Ingrédients" column (column C), shifting it and everything after
# it one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; existing C..G shift to D..H.
$ws.Columns.Item(3).Insert()

# New column width matches column B (19.44140625), since both now share it.
$ws.Columns.Item(3).ColumnWidth = 19.44140625

# Header for the new column.
$ws.Cells.Item(1, 3).Value = "Quantite"

# Numeric quantities parsed out of the existing "Quantité Ingrédients" text
# values (30g, 15cl, 1c.à.c, _, 4l, 2u -> 30, 15, 1, 0, 4, 2).
$ws.Cells.Item(2, 3).Value = 30
$ws.Cells.Item(3, 3).Value = 15
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 3).Value = 4
$ws.Cells.Item(9, 3).Value = 2

# Restore the active selection to what the author left it at.
$ws.Range("C9").Select()
